$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Japão"
$ws.Range("A15").Value = "Tonquim"
$ws.Range("A16").Value = "Suíça"
$ws.Range("A18").Value = "Alsácia"
$ws.Range("A19").Value = "Europa"
$ws.Range("A21").Value = "Coreia"
$ws.Range("A22").Value = "Dalmácia"
$ws.Range("A23").Value = "Lituânia"
$ws.Range("A25").Value = "Inglaterra"
$ws.Range("A26").Value = "México"
